$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.913.00'
$ws.Range("E2").Value = '  +0.11%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.533.03'
$ws.Range("E3").Value = '  +0.87%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.75'
$ws.Range("E5").Value = '  -0.38%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '196.62'
$ws.Range("E6").Value = '  +2.59%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.624'
$ws.Range("E7").Value = '  -0.22%  '

# Row 8
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$ws.Range("E9").Value = '  -4.27%  '

# Row 10
$ws.Range("E10").Value = '  -2.61%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.48'
$ws.Range("E11").Value = '  +0.02%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000303'
$ws.Range("E12").Value = '  -1.19%  '

# Row 13
$ws.Range("E13").Value = '  -1.29%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.088.52'
$ws.Range("E14").Value = '  +0.67%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '600.86'
$ws.Range("E15").Value = '  -3.41%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '12.79'
$ws.Range("E16").Value = '  +0.60%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '70.002.19'
$ws.Range("E17").Value = '  +0.13%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.03'
$ws.Range("E18").Value = '  +0.76%  '

# Row 19
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.530.20'
$ws.Range("E19").Value = '  +0.69%  '

# Row 20
$ws.Range("B20").Value = 'TRON'
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.123'
$ws.Range("E20").Value = '  +1.98%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.987'
$ws.Range("E21").Value = '  -0.50%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.00'
$ws.Range("E22").Value = '  +1.51%  '

# Row 23
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '103.53'
$ws.Range("E23").Value = '  -2.19%  '

# Row 24
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.17'
$ws.Range("E24").Value = '  +3.38%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.64'
$ws.Range("E25").Value = '  -0.18%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.06'
$ws.Range("E26").Value = '  +0.66%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.81'
$ws.Range("E27").Value = '  -1.54%  '

# Row 28
$ws.Range("E28").Value = '  -2.51%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.33'
$ws.Range("E29").Value = '  -2.41%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.10'
$ws.Range("E30").Value = '  +0.44%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.27'
$ws.Range("E31").Value = '  +3.28%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.36'
$ws.Range("E32").Value = '  -1.83%  '

# Row 33
$ws.Range("E33").Value = '  +0.25%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.77'
$ws.Range("E34").Value = '  -0.89%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.18'
$ws.Range("E35").Value = '  +2.89%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.774.70'
$ws.Range("E36").Value = '  +1.76%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0₃0816'
$ws.Range("E37").Value = '  +2.73%  '

# Row 38
$ws.Range("E38").Value = '  +0.15%  '

# Row 39
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.392'
$ws.Range("E39").Value = '  +0.29%  '

# Row 40
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '505.02'
$ws.Range("E40").Value = '  -3.26%  '

# Row 41
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.57'
$ws.Range("E41").Value = '  -0.09%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '36.48'
$ws.Range("E42").Value = '  -0.83%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.134'
$ws.Range("E43").Value = '  -3.13%  '

# Row 44
$ws.Range("E44").Value = '  -2.83%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.140'
$ws.Range("E45").Value = '  -0.61%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.82'
$ws.Range("E46").Value = '  -1.17%  '

# Row 47
$ws.Range("E47").Value = '  -2.67%  '

# Row 48
$ws.Range("E48").Value = '  +0.17%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.49'
$ws.Range("E49").Value = '  -2.91%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000250'
$ws.Range("E50").Value = '  +6.27%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.34'
$ws.Range("E51").Value = '  +3.28%  '

